$wb = $excel.ActiveWorkbook

# --- Remove the stray "Sheet1" tab (the 5x5 data that had been duplicated there) ---
$wb.Worksheets.Item("Sheet1").Delete()

# --- Work on the "3 x 3" sheet: add the partial magic-square algorithm scratch area ---
$ws = $wb.Worksheets.Item("3 x 3")

# Row 8: two scratch formula cells
$ws.Range("H8").Formula = "=1-1"
$ws.Range("H8").NumberFormat = "0"
$ws.Range("I8").Formula = "=0+1"

# Row 9: empty placeholder cell (reuses the "d-mmm" number format already used at C16)
$ws.Range("H9").NumberFormat = "d-mmm"
$ws.Range("H9").Value = ""

# Rows 17-21: row/column offsets worked out for each of the remaining numbers 5-9
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 2

$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2

$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0

$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 1

# Sheet view bookkeeping for "3 x 3": becomes the active tab, scrolled down one row,
# selection parked on D20
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("D20").Select()
$ws.Activate()

# "5 x 5" tab should no longer be the selected tab
$wb.Worksheets.Item("5 x 5").Select()
$ws.Activate()
